$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 100
$ws.Range("K9").Value = 100
$ws.Range("M9").Value = 69
$ws.Range("H17").Value = 1884.4
$ws.Range("J17").Value = 2263.4285
$ws.Range("L17").Value = 6790.2855
$ws.Range("N17").Value = -7126.2855
$ws.Range("H31").Value = 425
$ws.Range("I31").Value = 425
$ws.Range("K31").Value = 1275
$ws.Range("M31").Value = -1045
$ws.Range("H100").Value = 3433.6924
$ws.Range("I100").Value = 1483.6923
$ws.Range("J100").Value = 5383.6924
$ws.Range("K100").Value = 1483.6923
$ws.Range("L100").Value = 5383.6924
$ws.Range("M100").Value = -942.6922999999999
$ws.Range("N100").Value = -6465.6924
$ws.Range("H132").Value = 1707.7222
$ws.Range("I132").Value = 1141.6538
$ws.Range("K132").Value = 3424.9614
$ws.Range("M132").Value = -894.9614000000001
$ws.Range("H138").Value = 3049.742
$ws.Range("I138").Value = 2189
$ws.Range("J138").Value = 3300.7917
$ws.Range("K138").Value = 6567
$ws.Range("L138").Value = 9902.375100000001
$ws.Range("M138").Value = -1427
$ws.Range("N138").Value = -20182.3751
$ws.Range("H141").Value = 5034.4
$ws.Range("I141").Value = 5134.857
$ws.Range("K141").Value = 15404.571
$ws.Range("M141").Value = -10224.571

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3216.7942
$ws.Range("I102").Value = 1680.5927
$ws.Range("K102").Value = 1680.5927
$ws.Range("M102").Value = -58.59269999999992
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2544.0588
$ws.Range("I132").Value = 2608.4375
$ws.Range("K132").Value = 7825.3125
$ws.Range("M132").Value = -5295.3125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 950
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H58").Value = 37889.5
$ws.Range("J58").Value = 37889.5
$ws.Range("L58").Value = 37889.5
$ws.Range("N58").Value = -38477.5
$ws.Range("H81").Value = 27960
$ws.Range("J81").Value = 27960
$ws.Range("L81").Value = 27960
$ws.Range("N81").Value = -30082
$ws.Range("H84").Value = 27960
$ws.Range("J84").Value = 27960
$ws.Range("L84").Value = 83880
$ws.Range("N84").Value = -94488
$ws.Range("H99").Value = 2453.238
$ws.Range("J99").Value = 3835.9092
$ws.Range("L99").Value = 3835.9092
$ws.Range("N99").Value = -6831.9092
$ws.Range("H132").Value = 128871.375
$ws.Range("J132").Value = 128871.375
$ws.Range("L132").Value = 128871.375
$ws.Range("N132").Value = -138991.375
$ws.Range("H135").Value = 86570.71000000001
$ws.Range("J135").Value = 86570.71000000001
$ws.Range("L135").Value = 86570.71000000001
$ws.Range("N135").Value = -96710.71000000001
$ws.Range("H137").Value = 72186.336
$ws.Range("J137").Value = 72186.336
$ws.Range("L137").Value = 72186.336
$ws.Range("N137").Value = -82386.336
$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 299999.5
$ws.Range("L138").Value = 299999.5
$ws.Range("N138").Value = -310279.5
$ws.Range("H141").Value = 151000
$ws.Range("J141").Value = 151000
$ws.Range("L141").Value = 151000
$ws.Range("N141").Value = -161360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3156.2058
$ws.Range("I31").Value = 1446.909
$ws.Range("J31").Value = 6289.9165
$ws.Range("K31").Value = 1446.909
$ws.Range("L31").Value = 6289.9165
$ws.Range("M31").Value = -1151.909
$ws.Range("N31").Value = -6879.9165
$ws.Range("H32").Value = 4172.857
$ws.Range("I32").Value = 4172.857
$ws.Range("K32").Value = 4172.857
$ws.Range("M32").Value = -3856.857
$ws.Range("H34").Value = 3156.2058
$ws.Range("I34").Value = 1446.909
$ws.Range("J34").Value = 6289.9165
$ws.Range("K34").Value = 1446.909
$ws.Range("L34").Value = 6289.9165
$ws.Range("M34").Value = -1244.909
$ws.Range("N34").Value = -6693.9165
$ws.Range("H107").Value = 421.35
$ws.Range("I107").Value = 317
$ws.Range("J107").Value = 664.8333
$ws.Range("K107").Value = 317
$ws.Range("L107").Value = 664.8333
$ws.Range("M107").Value = 1603
$ws.Range("N107").Value = -4504.8333
$ws.Range("H134").Value = 1551.931
$ws.Range("I134").Value = 1518.1072
$ws.Range("K134").Value = 4554.321599999999
$ws.Range("M134").Value = -2019.321599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3499.75
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 3499.75
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10499.25
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10723.25
$ws.Range("H98").Value = 148.5
$ws.Range("J98").Value = 148.5
$ws.Range("L98").Value = 445.5
$ws.Range("N98").Value = -3441.5
$ws.Range("H123").Value = 932
$ws.Range("I123").Value = 932
$ws.Range("K123").Value = 2796
$ws.Range("M123").Value = -346
$ws.Range("H135").Value = 3499.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 3499.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 31497.75
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -36567.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H122").Value = 3135.2222
$ws.Range("I122").Value = 2309.5454
$ws.Range("K122").Value = 6928.6362
$ws.Range("M122").Value = -4478.6362
$ws.Range("H135").Value = 72500
$ws.Range("J135").Value = 72500
$ws.Range("L135").Value = 72500
$ws.Range("N135").Value = -82640

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 24666.666
$ws.Range("I40").Value = 24000
$ws.Range("K40").Value = 24000
$ws.Range("M40").Value = -23851
$ws.Range("H132").Value = 2645.658
$ws.Range("I132").Value = 2426.1562
$ws.Range("J132").Value = 3816.3333
$ws.Range("K132").Value = 7278.4686
$ws.Range("L132").Value = 11448.9999
$ws.Range("M132").Value = -4748.4686
$ws.Range("N132").Value = -16508.9999
